# Update: Threat Alert Report - 2026-01-30 06:29
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Market Threat Airline) width 25 -> 23 characters.
# ColumnWidth is expressed in "characters"; 22.14 round-trips to the
# stored OOXML width of 23 (same relationship the original 25 used: ~24.1-24.2).
$ws.Columns.Item(3).ColumnWidth = 22.14

# Row 2
# Force text format before assigning a date-looking string so the engine
# doesn't auto-convert it to a date serial number (keeps it literal text,
# matching the source report's plain "DD-MMM-YY" label column).
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "06-FEB-26"
$ws.Range("C2").Value = "Nile Air NP-118"
$ws.Range("D2").Value = 355
$ws.Range("E2").Value = 574
$ws.Range("F2").Value = -219
$ws.Range("G2").Value = 30
$ws.Range("I2").Value = 0

# Row 3
$ws.Range("D3").Value = 431
$ws.Range("E3").Value = 574
$ws.Range("F3").Value = -143

# Row 4
$ws.Range("C4").Value = "Nile Air NP-118"
$ws.Range("D4").Value = 355
$ws.Range("E4").Value = 574
$ws.Range("F4").Value = -219
$ws.Range("G4").Value = 30
$ws.Range("I4").Value = 0

# Row 5
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "13-FEB-26"
$ws.Range("D5").Value = 431
$ws.Range("E5").Value = 574
$ws.Range("F5").Value = -143

# Row 6
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "20-FEB-26"
$ws.Range("C6").Value = "EgyptAir MS-812"
$ws.Range("D6").Value = 383
$ws.Range("E6").Value = 536
$ws.Range("F6").Value = -153
$ws.Range("G6").Value = 46
$ws.Range("I6").Value = -16

# Row 7
$ws.Range("E7").Value = 536
$ws.Range("F7").Value = -20

# Row 8
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "06-MAR-26"
$ws.Range("D8").Value = 757
$ws.Range("E8").Value = 754
$ws.Range("F8").Value = 3
